$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace static values in column F (rows 2-72) with AVERAGE formulas
$ws.Range("F2").Formula = "=AVERAGE(B2:E2)"
$ws.Range("F3:F66").Formula = "=AVERAGE(B3:E3)"
$ws.Range("F67:F72").Formula = "=AVERAGE(B67:E67)"

# Make sure the whole column uses the percentage number format (0.00%)
$ws.Range("F2:F72").NumberFormat = "0.00%"

# Update the median formula to include row 72, and add average/max/min summary rows
$ws.Range("F75").Formula = "=MEDIAN(F2:F72)"

$ws.Range("E76").Value = "average"
$ws.Range("F76").Formula = "=AVERAGE(F2:F72)"

$ws.Range("E77").Value = "max"
$ws.Range("F77").Formula = "=MAX(F2:F72)"

$ws.Range("E78").Value = "min"
$ws.Range("F78").Formula = "=MIN(F2:F72)"

$ws.Range("F75:F78").NumberFormat = "0.00%"

# Update header F1: "ann" -> "ann  (ave quarters)"
$ws.Range("F1").Value = "ann  (ave quarters)"

# Widen column F to fit the new, longer header text
$ws.Columns.Item(6).ColumnWidth = 17

# Update the view: clear the old scroll/selection and select the new next-empty cell
$ws.Range("F79").Select() | Out-Null
